$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 9615.385
$ws.Range("I51").Value = 7500
$ws.Range("K51").Value = 7500
$ws.Range("M51").Value = -7016

$ws.Range("H98").Value = 2774.9167
$ws.Range("I98").Value = 1930.6
$ws.Range("K98").Value = 1930.6
$ws.Range("M98").Value = -432.5999999999999

$ws.Range("H122").Value = 2774.9167
$ws.Range("I122").Value = 1930.6
$ws.Range("K122").Value = 5791.799999999999
$ws.Range("M122").Value = -3341.799999999999

$ws.Range("H137").Value = 1000
$ws.Range("I137").Value = 1000
$ws.Range("K137").Value = 3000
$ws.Range("M137").Value = -450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2026.1538
$ws.Range("I134").Value = 2130.7273
$ws.Range("K134").Value = 6392.1819
$ws.Range("M134").Value = -3857.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 20733
$ws.Range("I25").Value = 30599.5
$ws.Range("K25").Value = 30599.5
$ws.Range("M25").Value = -30425.5

$ws.Range("H62").Value = 5581.6
$ws.Range("I62").Value = 7004.5
$ws.Range("J62").Value = 4633
$ws.Range("K62").Value = 7004.5
$ws.Range("L62").Value = 4633
$ws.Range("M62").Value = -6380.5
$ws.Range("N62").Value = -5881

$ws.Range("H65").Value = 5581.6
$ws.Range("I65").Value = 7004.5
$ws.Range("J65").Value = 4633
$ws.Range("K65").Value = 35022.5
$ws.Range("L65").Value = 23165
$ws.Range("M65").Value = -31902.5
$ws.Range("N65").Value = -29405

$ws.Range("H86").Value = 3972.5
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 3945
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 3945
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -6191

$ws.Range("H89").Value = 3972.5
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 3945
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 19725
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -30957

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 113.4
$ws.Range("J4").Value = 121
$ws.Range("L4").Value = 363
$ws.Range("N4").Value = -587

$ws.Range("H11").Value = 364.33334
$ws.Range("I11").Value = 399
$ws.Range("J11").Value = 347
$ws.Range("K11").Value = 1197
$ws.Range("L11").Value = 1041
$ws.Range("M11").Value = -1057
$ws.Range("N11").Value = -1321

$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 100
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 300
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -12

$ws.Range("H48").Value = 674.75
$ws.Range("J48").Value = 864.6667
$ws.Range("L48").Value = 2594.0001
$ws.Range("N48").Value = -3094.0001

$ws.Range("H52").Value = 500
$ws.Range("J52").Value = 500
$ws.Range("L52").Value = 1500
$ws.Range("N52").Value = -2032

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").ClearContents()
$ws.Range("N54").Value = 0

$ws.Range("H80").Value = 5000
$ws.Range("J80").Value = 5000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -16872

$ws.Range("H83").Value = 5000
$ws.Range("J83").Value = 5000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -54360

$ws.Range("H109").Value = 800
$ws.Range("I109").Value = 800
$ws.Range("K109").Value = 2400
$ws.Range("M109").Value = -1360

$ws.Range("H115").Value = 3413
$ws.Range("J115").Value = 4999.5
$ws.Range("L115").Value = 14998.5
$ws.Range("N115").Value = -17348.5

$ws.Range("H124").Value = 12400
$ws.Range("I124").Value = 4800
$ws.Range("J124").Value = 20000
$ws.Range("K124").Value = 14400
$ws.Range("L124").Value = 60000
$ws.Range("M124").Value = -9490
$ws.Range("N124").Value = -69820

$ws.Range("H125").Value = 20000
$ws.Range("J125").Value = 20000
$ws.Range("L125").Value = 60000
$ws.Range("N125").Value = -69840

$ws.Range("H126").Value = 12499.5
$ws.Range("J126").Value = 20000
$ws.Range("L126").Value = 60000
$ws.Range("N126").Value = -69880

$ws.Range("H129").Value = 20000
$ws.Range("J129").Value = 20000
$ws.Range("L129").Value = 60000
$ws.Range("N129").Value = -70000

$ws.Range("H130").Value = 11500
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 60000
$ws.Range("N130").Value = -70040

$ws.Range("H131").Value = 4449.9
$ws.Range("J131").Value = 3944.3333
$ws.Range("L131").Value = 11832.9999
$ws.Range("N131").Value = -21912.9999

$ws.Range("H140").Value = 2162.875
$ws.Range("I140").Value = 2162.875
$ws.Range("K140").Value = 6488.625
$ws.Range("M140").Value = -1308.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2157878
$ws.Range("I11").Value = 956827.0600000001
$ws.Range("J11").Value = 4800190
$ws.Range("K11").Value = 956827.0600000001
$ws.Range("L11").Value = 4800190
$ws.Range("M11").Value = -956688.0600000001
$ws.Range("N11").Value = -4800468

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 5000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 5000
$ws.Range("L20").ClearContents()
$ws.Range("N20").Value = 0
$ws.Range("M20").Value = -4774

$ws.Range("H46").Value = 450
$ws.Range("I46").Value = 450
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 450
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -262

$ws.Range("H61").Value = 5499.1665
$ws.Range("I61").Value = 4999.2
$ws.Range("K61").Value = 4999.2
$ws.Range("M61").Value = -4797.2

$ws.Range("H113").Value = 5499.1665
$ws.Range("I113").Value = 4999.2
$ws.Range("K113").Value = 4999.2
$ws.Range("M113").Value = -2829.2

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").ClearContents()
$ws.Range("N127").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 50017
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws.Range("H35").Value = 50017
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H62").Value = 2501
$ws.Range("I62").Value = 2002
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2002
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1378
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 2501
$ws.Range("I65").Value = 2002
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 10010
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -6890
$ws.Range("N65").Value = -21240
